$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# New values in columns E/F for rows 11-13
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 4

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 5

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 6

# Row 14 sums
$ws.Range("E14").Formula = "=SUM(E11:E12)"
$ws.Range("F14").Formula = "=SUM(F11:F12)"

# Update selection to match the new active cell
$ws.Activate()
$ws.Range("F14").Select()
